$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text formatting is preserved (avoid Excel auto-converting numeric-looking strings)
$cells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "D8", "E8", "D9", "E9", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "E22", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "E28", "E29", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D35", "E35", "D36", "E36", "D38", "E38", "E39", "E40", "E41", "D42", "E42", "D43", "E43", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "E49", "D50", "E50", "D51", "E51")
foreach ($c in $cells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.294.06"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "1.793.60"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "224.72"
$ws.Range("E5").Value = "  -1.86%  "
$ws.Range("D6").Value = "0.596"
$ws.Range("E6").Value = "  +2.73%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "36.21"
$ws.Range("E8").Value = "  +3.36%  "
$ws.Range("D9").Value = "0.291"
$ws.Range("E9").Value = "  -3.55%  "
$ws.Range("E10").Value = "  -3.83%  "
$ws.Range("D11").Value = "0.0962"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "2.052.16"
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").Value = "11.16"
$ws.Range("E13").Value = "  -2.41%  "
$ws.Range("D14").Value = "1.798.38"
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("D15").Value = "0.629"
$ws.Range("E15").Value = "  -2.88%  "
$ws.Range("D16").Value = "34.274.52"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "4.38"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "68.51"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").Value = "241.87"
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("D20").Value = "0.0₃0767"
$ws.Range("E20").Value = "  -4.48%  "
$ws.Range("D21").Value = "11.24"
$ws.Range("E21").Value = "  -3.39%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  -3.52%  "
$ws.Range("D24").Value = "2.18"
$ws.Range("E24").Value = "  +3.75%  "
$ws.Range("D25").Value = "170.71"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").Value = "7.95"
$ws.Range("E26").Value = "  +4.93%  "
$ws.Range("D27").Value = "17.24"
$ws.Range("E27").Value = "  +2.20%  "
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("D31").Value = "3.76"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").Value = "3.86"
$ws.Range("E32").Value = "  -3.31%  "
$ws.Range("D33").Value = "0.0510"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("D35").Value = "1.356.50"
$ws.Range("E35").Value = "  -3.29%  "
$ws.Range("D36").Value = "0.642"
$ws.Range("E36").Value = "  -5.85%  "
$ws.Range("D38").Value = "2.35"
$ws.Range("E38").Value = "  -7.80%  "
$ws.Range("E39").Value = "  -3.96%  "
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("E41").Value = "  -3.24%  "
$ws.Range("D42").Value = "80.30"
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("D43").Value = "0.931"
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("E44").Value = "  +5.11%  "
$ws.Range("D45").Value = "12.94"
$ws.Range("E45").Value = "  -6.90%  "
$ws.Range("D46").Value = "0.0494"
$ws.Range("E46").Value = "  -4.20%  "
$ws.Range("D47").Value = "1.954.19"
$ws.Range("E47").Value = "  -1.50%  "
$ws.Range("D48").Value = "5.73"
$ws.Range("E48").Value = "  -5.42%  "
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "101.80"
$ws.Range("E50").Value = "  -3.58%  "
$ws.Range("D51").Value = "0.0₆0118"
$ws.Range("E51").Value = "  -9.74%  "

# Restore original default cell style (no explicit number format) to match original formatting
foreach ($c in $cells) {
    $ws.Range($c).Style = "Normal"
}
